# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit refresh to the Leve profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW) per the diff of Coeurl_Profits.xlsx

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1785.7
$ws.Range("I86").Value = 1830
$ws.Range("J86").Value = 1766.7142
$ws.Range("K86").Value = 1830
$ws.Range("L86").Value = 1766.7142
$ws.Range("M86").Value = -707
$ws.Range("N86").Value = -4012.7142

$ws.Range("H89").Value = 1785.7
$ws.Range("I89").Value = 1830
$ws.Range("J89").Value = 1766.7142
$ws.Range("K89").Value = 9150
$ws.Range("L89").Value = 8833.571
$ws.Range("M89").Value = -3534
$ws.Range("N89").Value = -20065.571

$ws.Range("H103").Value = 391
$ws.Range("I103").Value = 875
$ws.Range("K103").Value = 2625
$ws.Range("M103").Value = -2039

$ws.Range("H138").Value = 6412983
$ws.Range("I138").Value = 1493.05
$ws.Range("J138").Value = 8623842
$ws.Range("K138").Value = 4479.15
$ws.Range("L138").Value = 25871526
$ws.Range("M138").Value = 660.8500000000004
$ws.Range("N138").Value = -25881806

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1164.95
$ws.Range("I2").Value = 1039.2941
$ws.Range("J2").Value = 1877
$ws.Range("K2").Value = 1039.2941
$ws.Range("L2").Value = 1877
$ws.Range("M2").Value = -926.2941000000001
$ws.Range("N2").Value = -2103

$ws.Range("H32").Value = 7128.662
$ws.Range("I32").Value = 3818.2834
$ws.Range("J32").Value = 25185.273
$ws.Range("K32").Value = 3818.2834
$ws.Range("L32").Value = 25185.273
$ws.Range("M32").Value = -3531.2834
$ws.Range("N32").Value = -25759.273

$ws.Range("H61").Value = 4570.619
$ws.Range("I61").Value = 3468.3
$ws.Range("K61").Value = 3468.3
$ws.Range("M61").Value = -3256.3

$ws.Range("H74").Value = 11995.363
$ws.Range("I74").Value = 2494.25
$ws.Range("K74").Value = 2494.25
$ws.Range("M74").Value = -1620.25

$ws.Range("H77").Value = 11995.363
$ws.Range("I77").Value = 2494.25
$ws.Range("K77").Value = 12471.25
$ws.Range("M77").Value = -8103.25

$ws.Range("H112").Value = 65166.668
$ws.Range("J112").Value = 65166.668
$ws.Range("L112").Value = 65166.668
$ws.Range("N112").Value = -68120.66800000001

$ws.Range("H116").Value = 1164.95
$ws.Range("I116").Value = 1039.2941
$ws.Range("J116").Value = 1877
$ws.Range("K116").Value = 1039.2941
$ws.Range("L116").Value = 1877
$ws.Range("M116").Value = 1254.7059
$ws.Range("N116").Value = -6465

$ws.Range("H132").Value = 2334.6487
$ws.Range("I132").Value = 2072.4666
$ws.Range("K132").Value = 6217.399800000001
$ws.Range("M132").Value = -3687.399800000001

$ws.Range("H136").Value = 4570.619
$ws.Range("I136").Value = 3468.3
$ws.Range("K136").Value = 10404.9
$ws.Range("M136").Value = -7854.900000000001

$ws.Range("H141").Value = 115000
$ws.Range("J141").Value = 115000
$ws.Range("L141").Value = 115000
$ws.Range("N141").Value = -125360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1164.95
$ws.Range("I3").Value = 1039.2941
$ws.Range("J3").Value = 1877
$ws.Range("K3").Value = 1039.2941
$ws.Range("L3").Value = 1877
$ws.Range("M3").Value = -925.2941000000001
$ws.Range("N3").Value = -2105

$ws.Range("H86").Value = 3737.1538
$ws.Range("I86").Value = 3750
$ws.Range("J86").Value = 3716.6
$ws.Range("K86").Value = 3750
$ws.Range("L86").Value = 3716.6
$ws.Range("M86").Value = -2627
$ws.Range("N86").Value = -5962.6

$ws.Range("H89").Value = 3737.1538
$ws.Range("I89").Value = 3750
$ws.Range("J89").Value = 3716.6
$ws.Range("K89").Value = 18750
$ws.Range("L89").Value = 18583
$ws.Range("M89").Value = -13134
$ws.Range("N89").Value = -29815

$ws.Range("H94").Value = 2715.1667
$ws.Range("I94").Value = 1768.0769
$ws.Range("J94").Value = 5177.6
$ws.Range("K94").Value = 1768.0769
$ws.Range("L94").Value = 5177.6
$ws.Range("M94").Value = -1317.0769
$ws.Range("N94").Value = -6079.6

$ws.Range("H105").Value = 2133.4707
$ws.Range("I105").Value = 1786.6666
$ws.Range("J105").Value = 4734.5
$ws.Range("K105").Value = 1786.6666
$ws.Range("L105").Value = 4734.5
$ws.Range("M105").Value = -39.66660000000002
$ws.Range("N105").Value = -8228.5

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 65535.25
$ws.Range("I31").Value = 73913.92999999999
$ws.Range("J31").Value = 6884.5
$ws.Range("K31").Value = 73913.92999999999
$ws.Range("L31").Value = 6884.5
$ws.Range("M31").Value = -73618.92999999999
$ws.Range("N31").Value = -7474.5

$ws.Range("H34").Value = 65535.25
$ws.Range("I34").Value = 73913.92999999999
$ws.Range("J34").Value = 6884.5
$ws.Range("K34").Value = 73913.92999999999
$ws.Range("L34").Value = 6884.5
$ws.Range("M34").Value = -73711.92999999999
$ws.Range("N34").Value = -7288.5

$ws.Range("H60").Value = 18128.857
$ws.Range("J60").Value = 18380.6
$ws.Range("L60").Value = 18380.6
$ws.Range("N60").Value = -19402.6

$ws.Range("H122").Value = 3196.5454
$ws.Range("I122").Value = 2312.7144
$ws.Range("J122").Value = 4743.25
$ws.Range("K122").Value = 6938.1432
$ws.Range("L122").Value = 14229.75
$ws.Range("M122").Value = -4488.1432
$ws.Range("N122").Value = -19129.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 96.40000000000001
$ws.Range("I2").Value = 96.28570999999999
$ws.Range("J2").Value = 98
$ws.Range("K2").Value = 577.71426
$ws.Range("L2").Value = 588
$ws.Range("M2").Value = -464.71426
$ws.Range("N2").Value = -814

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H108").Value = 4379.8
$ws.Range("I108").Value = 474.75
$ws.Range("K108").Value = 1424.25
$ws.Range("M108").Value = 1455.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws.Range("H97").Value = 6000
$ws.Range("I97").Value = 6000
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 6000
$ws.Range("L97").Value = 0
$ws.Range("N97").Value = -5504
$ws.Range("M97").ClearContents()

$ws.Range("H102").Value = 43480620
$ws.Range("I102").Value = 2162.5
$ws.Range("K102").Value = 2162.5
$ws.Range("M102").Value = -540.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3043.7646
$ws.Range("I22").Value = 1848.5
$ws.Range("J22").Value = 3411.5386
$ws.Range("K22").Value = 1848.5
$ws.Range("L22").Value = 3411.5386
$ws.Range("M22").Value = -1553.5
$ws.Range("N22").Value = -4001.5386

$ws.Range("H27").Value = 3043.7646
$ws.Range("I27").Value = 1848.5
$ws.Range("J27").Value = 3411.5386
$ws.Range("K27").Value = 1848.5
$ws.Range("L27").Value = 3411.5386
$ws.Range("M27").Value = -1741.5
$ws.Range("N27").Value = -3625.5386

$ws.Range("H46").Value = 1198.75
$ws.Range("J46").Value = 1800
$ws.Range("L46").Value = 1800
$ws.Range("N46").Value = -2176

$ws.Range("H93").Value = 1270.7368
$ws.Range("I93").Value = 1302.9375
$ws.Range("K93").Value = 1302.9375
$ws.Range("M93").Value = -54.9375

$ws.Range("H110").Value = 46333.332
$ws.Range("J110").Value = 46333.332
$ws.Range("L110").Value = 46333.332
$ws.Range("N110").Value = -54513.332

$ws.Range("H132").Value = 4765.643
$ws.Range("I132").Value = 4338.636
$ws.Range("K132").Value = 13015.908
$ws.Range("M132").Value = -10485.908

$ws.Range("H141").Value = 74333
$ws.Range("J141").Value = 74333
$ws.Range("L141").Value = 74333
$ws.Range("N141").Value = -84693
